# feat: custom captrow and datarow
# Insert a new header row ("descrow" marker row) above the existing
# header row on the "Item" sheet, shifting all existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item")
$ws.Activate()

# Insert a brand-new row 1; existing rows 1-4 become rows 2-5.
$ws.Rows.Item(1).Insert()

# Fill the new row with the literal marker text "descrow" across A1:X1.
$newRow = $ws.Range("A1:X1")
$newRow.Value = "descrow"

# Normalize style first (clears any fill/number-format inherited from the
# column definitions), then re-apply the centered alignment so every cell
# in the row shares one uniform style.
$newRow.Style = "Normal"
$newRow.HorizontalAlignment = -4108
$newRow.VerticalAlignment = -4108

# Update the visible selection to C8 (matches the saved view state).
$ws.Range("C8").Select()
